$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.574.08'
$ws.Range('E2').Value = '  +4.33%  '
$ws.Range('D3').Value = '1.598.75'
$ws.Range('E3').Value = '  +3.09%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.993'
$ws.Range('E4').Value = '  -0.68%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.57'
$ws.Range('E5').Value = '  +1.65%  '
$ws.Range('E6').Value = '  +6.54%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.992'
$ws.Range('E7').Value = '  -0.78%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '26.69'
$ws.Range('E8').Value = '  +11.93%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.250'
$ws.Range('E9').Value = '  +3.23%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0598'
$ws.Range('E10').Value = '  +2.80%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0911'
$ws.Range('E11').Value = '  +2.40%  '
$ws.Range('D12').Value = '1.823.23'
$ws.Range('E12').Value = '  +2.88%  '
$ws.Range('D13').Value = '1.586.24'
$ws.Range('E13').Value = '  +2.36%  '
$ws.Range('D14').Value = '29.533.08'
$ws.Range('E14').Value = '  +4.28%  '
$ws.Range('E15').Value = '  +3.93%  '
$ws.Range('E16').Value = '  +3.70%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.33'
$ws.Range('E17').Value = '  +4.20%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '243.01'
$ws.Range('E18').Value = '  +6.81%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.58'
$ws.Range('E19').Value = '  +3.07%  '
$ws.Range('D20').Value = '0.0₃0694'
$ws.Range('E20').Value = '  +2.74%  '
$ws.Range('E21').Value = '  -0.62%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.05'
$ws.Range('E22').Value = '  +3.70%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.29'
$ws.Range('E23').Value = '  +4.32%  '
$ws.Range('E24').Value = '  +3.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.89'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.29'
$ws.Range('E26').Value = '  +3.84%  '
$ws.Range('E27').Value = '  +5.51%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.41'
$ws.Range('E28').Value = '  +2.74%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.993'
$ws.Range('E29').Value = '  -0.73%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0473'
$ws.Range('E30').Value = '  +1.27%  '
$ws.Range('E31').Value = '  -0.08%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.24'
$ws.Range('E32').Value = '  +2.58%  '
$ws.Range('D33').Value = '1.434.08'
$ws.Range('E33').Value = '  +3.59%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.11'
$ws.Range('E34').Value = '  +3.52%  '
$ws.Range('E35').Value = '  -2.91%  '
$ws.Range('E36').Value = '  +2.96%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.80'
$ws.Range('E37').Value = '  +8.86%  '
$ws.Range('E38').Value = '  -1.67%  '
$ws.Range('E39').Value = '  +2.88%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.536'
$ws.Range('E40').Value = '  +5.37%  '
$ws.Range('E41').Value = '  +3.14%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '54.37'
$ws.Range('E42').Value = '  +28.84%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.805'
$ws.Range('E43').Value = '  +3.65%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.993'
$ws.Range('E44').Value = '  -0.67%  '
$ws.Range('E45').Value = '  +3.56%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '65.52'
$ws.Range('E46').Value = '  +5.75%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.37'
$ws.Range('D48').Value = '1.733.90'
$ws.Range('E48').Value = '  +2.93%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '86.23'
$ws.Range('E49').Value = '  +0.91%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.835'
$ws.Range('E50').Value = '  -3.59%  '
$ws.Range('E51').Value = '  +1.75%  '
